$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing status cells (Column D) ---
$ws.Range("D9").Value = "Removed"
$ws.Range("D16").Value = "Done"
$ws.Range("D17").Value = "Removed"

# --- Add new rows 20-25, copying the formatting of row 19 (the last
#     existing "Ambience" row) so the new rows keep the same styles ---
$ws.Range("A19:D19").Copy()
$ws.Range("A20:D25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A20").Value = "Ambience"
$ws.Range("B20").Value = "Warnings"
$ws.Range("C20").Value = "Alert sound on the top of the screen"
$ws.Range("D20").Value = "Added"

$ws.Range("A21").Value = "Ambience"
$ws.Range("B21").Value = "Engine hovering"
$ws.Range("C21").Value = "Sci-fi engine hovering sound scattering"
$ws.Range("D21").Value = "Added"

$ws.Range("A22").Value = "Ambience"
$ws.Range("B22").Value = "Future city surroundings"
$ws.Range("C22").Value = "Low pitch sci-fi digital sounds scattering"
$ws.Range("D22").Value = "Added"

$ws.Range("A23").Value = "Sound effect"
$ws.Range("B23").Value = "Freeze gun pickup"
$ws.Range("C23").Value = "Gun reloading sound, faster"
$ws.Range("D23").Value = "Added"

$ws.Range("A24").Value = "Sound effect"
$ws.Range("B24").Value = "Gravity gun pickup"
$ws.Range("C24").Value = "Weapon picking up sound"
$ws.Range("D24").Value = "Added"

$ws.Range("A25").Value = "Sound effect"
$ws.Range("B25").Value = "Sonic gun pickup"
$ws.Range("C25").Value = "Gun reloading sound, slower"
$ws.Range("D25").Value = "Added"

# --- Match the final selection left by the author ---
[void]$ws.Range("B21").Select()
